# Update "想去人数" (number of people interested) counts in column F
# across the "展览", "演出", and "全部类型" sheets, per the upstream
# data refresh recorded in the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3591
$ws1.Range("F5").Value = 3591
$ws1.Range("F6").Value = 262
$ws1.Range("F7").Value = 5119
$ws1.Range("F8").Value = 534
$ws1.Range("F9").Value = 362
$ws1.Range("F13").Value = 96
$ws1.Range("F14").Value = 35
$ws1.Range("F15").Value = 702
$ws1.Range("F22").Value = 4920
$ws1.Range("F25").Value = 12
$ws1.Range("F26").Value = 6050
$ws1.Range("F30").Value = 345
$ws1.Range("F36").Value = 1029
$ws1.Range("F41").Value = 1006

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 26

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3591
$ws4.Range("F8").Value = 3591
$ws4.Range("F9").Value = 262
$ws4.Range("F10").Value = 5119
$ws4.Range("F11").Value = 534
$ws4.Range("F12").Value = 362
$ws4.Range("F16").Value = 96
$ws4.Range("F17").Value = 35
$ws4.Range("F18").Value = 702
$ws4.Range("F26").Value = 4920
$ws4.Range("F29").Value = 12
$ws4.Range("F30").Value = 6050
$ws4.Range("F34").Value = 345
$ws4.Range("F38").Value = 26
$ws4.Range("F41").Value = 1029
$ws4.Range("F46").Value = 1006
